$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11, shifting existing rows 11-25 down to 12-26.
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with data.
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(11, 3).Value = "La Araucanía"
$ws.Cells.Item(11, 4).Value = 44482
$ws.Cells.Item(11, 4).NumberFormat = $ws.Cells.Item(12, 4).NumberFormat
$ws.Cells.Item(11, 5).Value = 9
$ws.Cells.Item(11, 6).Value = 300000000
$ws.Cells.Item(11, 7).Value = "Espárragos"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 120
$ws.Cells.Item(11, 11).Value = 1500
$ws.Cells.Item(11, 12).Value = 1500
$ws.Cells.Item(11, 13).Value = 1500
$ws.Cells.Item(11, 14).Value = "$/kilo"
$ws.Cells.Item(11, 15).Value = "Región del Maule"
$ws.Cells.Item(11, 16).Value = 1500
$ws.Cells.Item(11, 17).Value = 1
$ws.Cells.Item(11, 18).Value = "Hortaliza"
